$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header / account holder details
$ws.Range("C2").Value = "Hartmut"

# B3 holds a long card-number string that must stay text (inlineStr), not be
# auto-converted to a number. Stash B3's current formatting in a scratch cell,
# force Text format + assign the new value, then paste the original formatting
# back over it (keeping the value's text type) and wipe the scratch cell.
$ws.Range("B3").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("Z1").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("Z1").Clear()

$ws.Range("C3").Value = "Mohaupt"

# Opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 26.03.2024"

# Row 6
$ws.Range("B6").Value = "28.03."
$ws.Range("C6").Value = "29.03."
$ws.Range("D6").Value = "MCDONALDS Gräfenhainichen"
$ws.Range("E6").Value = "26,68-"

# Row 7
$ws.Range("B7").Value = "01.04."
$ws.Range("C7").Value = "02.04."
$ws.Range("D7").Value = "BEITRAG Allianz SE K-8042898"
$ws.Range("E7").Value = "54,32-"

# Row 8
$ws.Range("B8").Value = "02.04."
$ws.Range("C8").Value = "03.04."
$ws.Range("D8").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E8").Value = "24,62-"

# Row 9 - previously blank, now a new transaction row.
# Match the formatting used by the other transaction rows (B8/C8/D8 = style 8, E8 = style 17)
# by copying the row's formats across before filling in the new values.
$ws.Range("B8:E8").Copy()
$ws.Range("B9:E9").PasteSpecial(-4122)

$ws.Range("B9").Value = "05.04."
$ws.Range("C9").Value = "06.04."
$ws.Range("D9").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 84255032"
$ws.Range("E9").Value = "86,77-"

# Closing balance line
$ws.Range("D12").Value = "KONTOSTAND AM 10.04.2024"
$ws.Range("E12").Value = "192,39-"

# Next statement date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 17.04.2024"
